$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    # Fix the Url column (C): strip the trailing ") " left over from a
    # copy/paste bug so links resolve correctly.
    $urlCell = $ws.Cells.Item($r, 3)
    $url = $urlCell.Value2
    if ($url -ne $null -and $url.EndsWith(") ")) {
        $urlCell.Value2 = $url.Substring(0, $url.Length - 2)
    }

    # Normalize the Holder column (D) casing now that a bank.txt can list
    # multiple banks/holders.
    $holderCell = $ws.Cells.Item($r, 4)
    $holder = $holderCell.Value2
    if ($holder -eq "peanuts") {
        $holderCell.Value2 = "Peanuts"
    }
}
